$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new columns before column D ---
# This shifts the previous D:K quarterly data right by two columns (-> F:M)
# and opens up D:E for the two new quarters being added.
$ws.Range("D1:E1").EntireColumn.Insert()

# --- 2. Copy number/date formatting from column F into new columns D:E ---
# Done in three contiguous row blocks so the header-only rows 37 and 79
# (which have no D:K cells at all) are left untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Helper to bulk-assign a single worksheet row via a 2-D COM array ---
function Set-RowValues {
    param($WorksheetObj, $RangeAddress, $Values)
    $arr = New-Object "object[,]" 1, $Values.Count
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $WorksheetObj.Range($RangeAddress).Value = $arr
}

# --- 4. Write the refreshed financial data across columns D:M ---
Set-RowValues $ws "D7:M7" @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws "D8:M8" @(1122000, 1388400, 1086200, 998900, 916100, 1172400, 952200, 879100, 736900, 810300)
Set-RowValues $ws "D9:M9" @(240400, 295500, 218000, 184600, 157500, 193500, 166800, 172900, 163800, 179500)
Set-RowValues $ws "D10:M10" @(881600, 1092900, 868200, 814300, 758700, 978900, 785300, 706200, 573100, 630800)
Set-RowValues $ws "D12:M12" @(403400, 369700, 334100, 320600, 307800, 324300, 302300, 285500, 250500, 268200)
Set-RowValues $ws "D13:M13" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D14:M14" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D15:M15" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D17:M17" @(1150000, 1168700, 978800, 911400, 871200, 970700, 856400, 824700, 706800, 745200)
Set-RowValues $ws "D18:M18" @(-28000, 219600, 107400, 87600, 45000, 201700, 95800, 54400, 30200, 65100)
Set-RowValues $ws "D20:M20" @(-78500, -312300, 382900, 130200, 99900, 74900, 95900, 6100, -23300, 10400)
Set-RowValues $ws "D21:M21" @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws "D22:M22" @(62600, 58300, 55100, 47800, 48100, 47300, 56900, 37800, 31900, 27800)
Set-RowValues $ws "D23:M23" @(-169200, -150900, 435300, 169900, 96800, 229300, 134800, 22700, -25100, 47600)
Set-RowValues $ws "D24:M24" @(5200, 38100, 47800, 26600, 35300, 46500, 78500, 20200, 16000, 32200)
Set-RowValues $ws "D25:M25" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D26:M26" @(-174400, -189100, 387500, 143400, 61400, 182800, 56300, 2500, -41100, 15500)
Set-RowValues $ws "D27:M27" @(-176600, -169000, 353400, 157300, 51900, 182400, 48500, 7600, 93800, 3500)
Set-RowValues $ws "D28:M28" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D29:M29" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D30:M30" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D31:M31" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D32:M32" @(78500, 312300, -382900, -130200, -99900, -74900, -95900, -6100, 23300, -10400)
Set-RowValues $ws "D33:M33" @(-176600, -169000, 353400, 157300, 51900, 182400, 48500, 7600, 93800, 3500)
Set-RowValues $ws "D34:M34" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D35:M35" @(-176600, -169000, 353400, 157300, 51900, 182400, 48500, 7600, 93800, 3500)
Set-RowValues $ws "D38:M38" @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws "D41:M41" @(3195300, 3189300, 3271100, 2560500, 2707400, 2719400, 2403700, 2751500, 2681000, 4475200)
Set-RowValues $ws "D42:M42" @(5454500, 5608900, 4923900, 5064200, 4174800, 4016400, 3600300, 2166800, 2052400, 901200)
Set-RowValues $ws "D43:M43" @(1084900, 1008700, 929200, 768500, 760300, 869000, 781600, 700300, 787100, 638400)
Set-RowValues $ws "D44:M44" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D45:M45" @(2048200, 2728400, 1901600, 1448600, 1231200, 1780900, 1345100, 1160800, 1156400, 1579800)
Set-RowValues $ws "D46:M46" @(11782900, 12535300, 11025800, 9841800, 8818200, 9385600, 8130800, 6779500, 6676900, 7594700)
Set-RowValues $ws "D47:M47" @(4022400, 3806000, 4041200, 3777600, 3865800, 3403800, 3353200, 3242000, 3104700, 2795900)
Set-RowValues $ws "D48:M48" @(871500, 866300, 864000, 858600, 833500, 835800, 838100, 809300, 813200, 810700)
Set-RowValues $ws "D49:M49" @(10662200, 10681100, 10629100, 10433100, 10402500, 10411500, 10396600, 10171700, 10185800, 8350900)
Set-RowValues $ws "D50:M50" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D51:M51" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D52:M52" @(240100, 215800, 248600, 197100, 193200, 187900, 237600, 221700, 221400, 162100)
Set-RowValues $ws "D53:M53" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D54:M54" @(27579000, 28104400, 26808800, 25108200, 24078000, 24224600, 22956300, 21224300, 21002100, 19714200)
Set-RowValues $ws "D57:M57" @(1738500, 2018200, 1648500, 1397400, 1107000, 1233900, 1043200, 933200, 1058600, 1116900)
Set-RowValues $ws "D58:M58" @(5344400, 5346600, 3630400, 3053700, 2421500, 2202200, 1571200, 1555000, 1001600, 1973700)
Set-RowValues $ws "D59:M59" @(3125400, 2885100, 2918200, 2564100, 2790900, 2793000, 2545300, 2126200, 2345600, 1733900)
Set-RowValues $ws "D60:M60" @(10208200, 10249900, 8197100, 7015200, 6257300, 6229100, 5159800, 4614400, 4405800, 4824500)
Set-RowValues $ws "D61:M61" @(3583500, 3867900, 4439200, 4314600, 4336500, 4986200, 5282800, 4988900, 5039200, 4074700)
Set-RowValues $ws "D62:M62" @(618400, 631600, 672100, 612300, 629700, 592500, 595800, 567500, 574100, 452200)
Set-RowValues $ws "D63:M63" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D64:M64" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D65:M65" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D66:M66" @(14709700, 15080500, 13590900, 12220400, 11487500, 12095000, 11334500, 10435400, 10598500, 9886300)
Set-RowValues $ws "D68:M68" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D69:M69" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D70:M70" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D71:M71" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D72:M72" @(2437900, "NA", "NA", 2460800, 2303500, 1272700, 1090300, 1020800, 1008900, 915000)
Set-RowValues $ws "D73:M73" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D74:M74" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D75:M75" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D76:M76" @(12869400, 13023900, 13217800, 12887800, 12590500, 12129600, 11621800, 10788900, 10403600, 9827900)
Set-RowValues $ws "D77:M77" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D80:M80" @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws "D81:M81" @(-176600, -169000, 353400, 157300, 51900, 182400, 48500, 7600, 93800, 3500)
Set-RowValues $ws "D83:M83" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D84:M84" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D85:M85" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D86:M86" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D87:M87" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D88:M88" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D89:M89" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D91:M91" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D92:M92" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D93:M93" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D94:M94" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D96:M96" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D97:M97" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D98:M98" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D99:M99" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D100:M100" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D101:M101" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D102:M102" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
